$wb = $excel.ActiveWorkbook

# The three summary worksheets each hold one statsmodels OLS text dump in B2.
# Only the "Time:" line (when the summary was generated) changed between saves;
# everything else in each block stays the same.
$ws1 = $wb.Worksheets.Item(1)   # sheet "5" -> Time: 20:51:41 -> 20:59:42
$ws2 = $wb.Worksheets.Item(2)   # sheet "4" -> Time: 20:51:41 -> 20:59:42
$ws3 = $wb.Worksheets.Item(3)   # sheet "3" -> Time: 20:51:41 -> 20:59:43

$updates = @(
    @{ Sheet = $ws1; Old = "20:51:41"; New = "20:59:42" },
    @{ Sheet = $ws2; Old = "20:51:41"; New = "20:59:42" },
    @{ Sheet = $ws3; Old = "20:51:41"; New = "20:59:43" }
)

foreach ($u in $updates) {
    $cell = $u.Sheet.Cells.Item(2, 2)
    $text = $cell.Value2
    $oldLine = "Time:                        " + $u.Old
    $newLine = "Time:                        " + $u.New
    $cell.Value2 = $text.Replace($oldLine, $newLine)
}
